$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Bump the cached "datetimeFigureOut" date field text from 11/23/2019 to
#    11/24/2019 on every slide layout and on the slide master (the Date
#    Placeholder shape on each).
# ---------------------------------------------------------------------------
$oldDate = "11/23/2019"
$newDate = "11/24/2019"

$master = $p.SlideMaster

# Slide master's own Date Placeholder shape.
for ($mi = 1; $mi -le $master.Shapes.Count; $mi++) {
    $mShape = $master.Shapes.Item($mi)
    if ($mShape.Name -like "Date Placeholder*") {
        if ($mShape.TextFrame.TextRange.Text -eq $oldDate) {
            $mShape.TextFrame.TextRange.Text = $newDate
        }
    }
}

# Each slide layout's Date Placeholder shape.
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    for ($si = 1; $si -le $layout.Shapes.Count; $si++) {
        $lShape = $layout.Shapes.Item($si)
        if ($lShape.Name -like "Date Placeholder*") {
            if ($lShape.TextFrame.TextRange.Text -eq $oldDate) {
                $lShape.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# ---------------------------------------------------------------------------
# 2) On slide 1's title shape, split the GitHub hyperlink run into two runs
#    ("https://github.com/bobk" + "/jirapresentations") so it matches the
#    run-split left behind by the author's edit, keeping the hyperlink on
#    both pieces.
# ---------------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
$titleShape = $slide1.Shapes.Item(1)
$titleRange = $titleShape.TextFrame.TextRange

$fullTitleText = $titleRange.Text
$url = "https://github.com/bobk/jirapresentations"
$splitAt = "https://github.com/bobk"

$urlIndex = $fullTitleText.IndexOf($url)
if ($urlIndex -ge 0) {
    $firstPart = $titleRange.Characters($urlIndex + 1, $splitAt.Length)
    $firstPart.Text = $splitAt
}
